# Fix du fix issue#70
$wb = $excel.ActiveWorkbook

# 1) Update the Date value on the Metadata sheet (row 8: "Date" / value)
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-02-16T15:11:19+00:00"

# 2) Insert a new "GoNogo -> MedicationRequest.status" mapping row into
#    "Mapping Table 10" right before the "Motif_attente" row.
$ws = $wb.Worksheets.Item("Mapping Table 10")
$ws.Rows.Item(8).Insert()

# Copy the formatting of the row above onto the newly inserted (blank) row
# so it keeps the same bordered-cell style as the rest of the table.
$ws.Range("A7:E7").Copy()
$ws.Range("A8:E8").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A8").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("C8").Value = "equivalent"
$ws.Range("D8").Value = "MedicationRequest.status"
$ws.Range("E8").ClearContents()
